{"js": "// Replace the multiplication problems in the table with the new set of\n// problems. Each old value is unique in the document, so a search +\n// replace keyed on the full old cell text (\"123\u00d74=\") is unambiguous and\n// keeps the surrounding run formatting (font/size) untouched.\nconst replacements = [\n  [\"276\u00d77=\", \"715\u00d76=\"],\n  [\"622\u00d73=\", \"138\u00d79=\"],\n  [\"656\u00d75=\", \"667\u00d72=\"],\n  [\"413\u00d78=\", \"713\u00d79=\"],\n  [\"324\u00d78=\", \"662\u00d73=\"],\n  [\"543\u00d72=\", \"676\u00d75=\"],\n  [\"733\u00d78=\", \"206\u00d74=\"],\n  [\"830\u00d75=\", \"807\u00d75=\"],\n  [\"216\u00d77=\", \"529\u00d74=\"],\n  [\"724\u00d79=\", \"259\u00d73=\"],\n  [\"422\u00d75=\", \"104\u00d74=\"],\n  [\"447\u00d78=\", \"194\u00d79=\"],\n  [\"175\u00d73=\", \"359\u00d73=\"],\n  [\"272\u00d78=\", \"368\u00d75=\"],\n  [\"363\u00d74=\", \"961\u00d75=\"],\n  [\"616\u00d73=\", \"115\u00d76=\"],\n  [\"692\u00d73=\", \"587\u00d74=\"],\n  [\"335\u00d73=\", \"176\u00d76=\"],\n  [\"386\u00d77=\", \"783\u00d79=\"],\n  [\"643\u00d73=\", \"448\u00d73=\"],\n  [\"741\u00d79=\", \"625\u00d72=\"],\n  [\"990\u00d73=\", \"839\u00d77=\"],\n  [\"663\u00d76=\", \"534\u00d75=\"],\n  [\"507\u00d74=\", \"981\u00d77=\"],\n  [\"714\u00d77=\", \"469\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication problems in the table with the new set of\n# problems. Each old value is unique in the document, so Find/Replace\n# keyed on the full old cell text (\"123\u00d74=\") is unambiguous and keeps\n# the surrounding run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n    \"276\u00d77=\" = \"715\u00d76=\"\n    \"622\u00d73=\" = \"138\u00d79=\"\n    \"656\u00d75=\" = \"667\u00d72=\"\n    \"413\u00d78=\" = \"713\u00d79=\"\n    \"324\u00d78=\" = \"662\u00d73=\"\n    \"543\u00d72=\" = \"676\u00d75=\"\n    \"733\u00d78=\" = \"206\u00d74=\"\n    \"830\u00d75=\" = \"807\u00d75=\"\n    \"216\u00d77=\" = \"529\u00d74=\"\n    \"724\u00d79=\" = \"259\u00d73=\"\n    \"422\u00d75=\" = \"104\u00d74=\"\n    \"447\u00d78=\" = \"194\u00d79=\"\n    \"175\u00d73=\" = \"359\u00d73=\"\n    \"272\u00d78=\" = \"368\u00d75=\"\n    \"363\u00d74=\" = \"961\u00d75=\"\n    \"616\u00d73=\" = \"115\u00d76=\"\n    \"692\u00d73=\" = \"587\u00d74=\"\n    \"335\u00d73=\" = \"176\u00d76=\"\n    \"386\u00d77=\" = \"783\u00d79=\"\n    \"643\u00d73=\" = \"448\u00d73=\"\n    \"741\u00d79=\" = \"625\u00d72=\"\n    \"990\u00d73=\" = \"839\u00d77=\"\n    \"663\u00d76=\" = \"534\u00d75=\"\n    \"507\u00d74=\" = \"981\u00d77=\"\n    \"714\u00d77=\" = \"469\u00d79=\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
